$qCases = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
MATCH (samp:sample)-->(c)
WHERE  samp.specific_sample_pathology in ['Pulmonary Carcinoma']
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
  coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
order by c.case_id asc
limit 100
'@

$qCasesStat = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma]
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
   count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$qSamples = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE  samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed,
        coalesce(diag.disease_term,'') AS Diagnosis, 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
order by samp.sample_id asc
limit 200
'@

$qStatShared = @'
MATCH (p:program)<--(s:study)<--(c)
MATCH (cf)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
MATCH (cf:file)-[*]->(c:case)
OPTIONAL MATCH (sf:file)-->(s)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct cf) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$qFiles = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
MATCH (f)-->(samp:sample)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
 MATCH (f)-[*]->(samp:sample)
WITH
        DISTINCT f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
        order by f.file_name asc
        limit 200
'@

$qStudyFiles = @'
MATCH (f:file)-->(s:study)
MATCH (s)<--(c:case)<--(diag:diagnosis)
MATCH (c)<--(demo:demographic)
MATCH (samp:sample)-->(c)
WHERE samp.specific_sample_pathology IN ['Pulmonary Carcinoma']
WITH
        DISTINCT f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
  order by 'File Name' asc
  limit 100
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Content is entered in the same order the author produced it in (this
# drives the order new entries land in xl/sharedStrings.xml), then
# formatting is applied as a separate pass.

# 1) New tab label for the duplicated/added row
$ws.Range("A5").Value = "StudyFilesTab"

# 2) CasesTab (row 2) query + stat query
$ws.Range("B2").Value = $qCases
$ws.Range("C2").Value = $qCasesStat

# 3) Shared stat query first lands on SamplesTab's stat cell
$ws.Range("C3").Value = $qStatShared

# 4) FilesTab query
$ws.Range("B4").Value = $qFiles

# 5) SamplesTab query
$ws.Range("B3").Value = $qSamples

# 6) StudyFilesTab query
$ws.Range("B5").Value = $qStudyFiles

# 7) Remaining cells reuse already-registered shared strings
$ws.Range("C4").Value = $qStatShared
$ws.Range("C5").Value = $qStatShared
$ws.Range("D5").Value = "TC09_Canine_Filter_SamplePatho-PulmoCarcinoma_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC09_Canine_Filter_SamplePatho-PulmoCarcinoma_WebData.xlsx"

# --- Formatting pass ---
$ws.Range("B2").WrapText = $true

$ws.Range("C2").WrapText = $true
$ws.Range("C2").VerticalAlignment = -4108

$ws.Range("B3").WrapText = $true
$ws.Range("B3").Font.Size = 18

$ws.Range("C3").WrapText = $true
$ws.Range("C3").VerticalAlignment = -4108

$ws.Range("B4").WrapText = $true
$ws.Range("B4").Font.Size = 18

$ws.Range("C4").WrapText = $true
$ws.Range("C4").VerticalAlignment = -4108

$ws.Range("B5").WrapText = $true
$ws.Range("B5").Font.Size = 18

$ws.Range("C5").WrapText = $true
$ws.Range("C5").Font.Size = 18

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 304.5
$ws.Rows.Item(3).RowHeight = 409.5
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Rows.Item(5).RowHeight = 409.5

# --- View state ---
$excel.ActiveWindow.Zoom = 69
$ws.Range("B5").Select()

Write-Host "edit complete"
